$wb = $excel.ActiveWorkbook

# Sheet "展览": F9 124 -> 126
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 126

# Sheet "演出": F3 29 -> 30
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 30

# Sheet "全部类型": F10 124 -> 126, F11 29 -> 30
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 126
$ws4.Range("F11").Value = 30
